# The deck's slide master (and therefore every ordinary slide) was using
# the "Integral" design theme while the notes master carried the default
# "Office Theme" palette. The edit flips which palette is "active": the
# design applied to the slides becomes the plain Office color palette
# (the notes-side theme keeps its own look).
#
# PowerPoint exposes the twelve theme colour slots for the design that is
# actually applied to the slides through Slide.ThemeColorScheme - each
# slot is an msoThemeColorSchemeIndex entry (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) with a settable .RGB. Driving all twelve of
# them to the stock Office values reproduces the "switch the active
# design to the Office Theme palette" edit.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in msoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1      - 000000
    (RGBVal 0xFF 0xFF 0xFF),  # lt1      - FFFFFF
    (RGBVal 0x44 0x54 0x6A),  # dk2      - 44546A
    (RGBVal 0xE7 0xE6 0xE6),  # lt2      - E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),  # accent1  - 5B9BD5
    (RGBVal 0xED 0x7D 0x31),  # accent2  - ED7D31
    (RGBVal 0xA5 0xA5 0xA5),  # accent3  - A5A5A5
    (RGBVal 0xFF 0xC0 0x00),  # accent4  - FFC000
    (RGBVal 0x44 0x72 0xC4),  # accent5  - 4472C4
    (RGBVal 0x70 0xAD 0x47),  # accent6  - 70AD47
    (RGBVal 0x05 0x63 0xC1),  # hlink    - 0563C1
    (RGBVal 0x95 0x4F 0x72)   # folHlink - 954F72
)

$p = $ppt.ActivePresentation

# Every slide shares the same slide master/design, so touching slide 1's
# theme colour scheme retints the one design used throughout the deck.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
